# Apply crypto price/volume updates (and the EnergySwap/Algorand row swap)
# described by the commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.580.18"
$ws.Range("E2").Value = "  +0.01%  "

$ws.Range("D3").Value = "1.922.46"
$ws.Range("E3").Value = "  -0.06%  "

$ws.Range("D4").Value = "'0.9983"
$ws.Range("E4").Value = "  -0.21%  "

$ws.Range("D5").Value = "'245.70"
$ws.Range("E5").Value = "  -0.65%  "

$ws.Range("D6").Value = "'0.9996"
$ws.Range("E6").Value = "  -0.10%  "

$ws.Range("D7").Value = "'0.4895"
$ws.Range("E7").Value = "  +3.26%  "

$ws.Range("D8").Value = "'0.2906"
$ws.Range("E8").Value = "  -0.16%  "

$ws.Range("D9").Value = "'0.06727"
$ws.Range("E9").Value = "  -0.98%  "

$ws.Range("D10").Value = "'111.22"
$ws.Range("E10").Value = "  +5.53%  "

$ws.Range("D11").Value = "'19.15"
$ws.Range("E11").Value = "  +3.79%  "

$ws.Range("D12").Value = "1.913.67"
$ws.Range("E12").Value = "  -0.49%  "

$ws.Range("D13").Value = "'0.07585"
$ws.Range("E13").Value = "  -1.80%  "

$ws.Range("D14").Value = "'5.294"
$ws.Range("E14").Value = "  -0.49%  "

$ws.Range("D15").Value = "'0.6717"
$ws.Range("E15").Value = "  -0.01%  "

$ws.Range("D16").Value = "'296.24"
$ws.Range("E16").Value = "  +3.05%  "

$ws.Range("D17").Value = "30.567.14"
$ws.Range("E17").Value = "  -0.16%  "

$ws.Range("D18").Value = "'13.03"
$ws.Range("E18").Value = "  +0.56%  "

$ws.Range("E19").Value = "  +0.03%  "

$ws.Range("D20").Value = "'0.000007589"
$ws.Range("E20").Value = "  -0.56%  "

$ws.Range("D21").Value = "'5.556"
$ws.Range("E21").Value = "  +1.93%  "

$ws.Range("D22").Value = "2.167.19"
$ws.Range("E22").Value = "  +0.15%  "

$ws.Range("D23").Value = "'1.0000"
$ws.Range("E23").Value = "  +0.00%  "

$ws.Range("D24").Value = "'6.471"
$ws.Range("E24").Value = "  +2.46%  "

$ws.Range("D25").Value = "'9.488"
$ws.Range("E25").Value = "  +0.81%  "

$ws.Range("D26").Value = "'164.89"
$ws.Range("E26").Value = "  -1.94%  "

$ws.Range("D27").Value = "'20.33"
$ws.Range("E27").Value = "  -2.56%  "

$ws.Range("D28").Value = "'2.104"
$ws.Range("E28").Value = "  -2.06%  "

$ws.Range("D29").Value = "'0.1078"
$ws.Range("E29").Value = "  -0.34%  "

$ws.Range("D30").Value = "'1.447"
$ws.Range("E30").Value = "  +5.83%  "

$ws.Range("D31").Value = "'4.163"
$ws.Range("E31").Value = "  -0.90%  "

$ws.Range("D32").Value = "'4.066"
$ws.Range("E32").Value = "  -1.42%  "

$ws.Range("D33").Value = "'0.05027"
$ws.Range("E33").Value = "  -0.40%  "

$ws.Range("D34").Value = "'0.7417"
$ws.Range("E34").Value = "  -0.11%  "

$ws.Range("D35").Value = "'1.142"
$ws.Range("E35").Value = "  -1.61%  "

$ws.Range("D36").Value = "'0.9998"
$ws.Range("E36").Value = "  -0.02%  "

$ws.Range("D37").Value = "'2.711"
$ws.Range("E37").Value = "  -1.26%  "

$ws.Range("D38").Value = "'0.02029"
$ws.Range("E38").Value = "  -2.29%  "

$ws.Range("D39").Value = "'2.679"
$ws.Range("E39").Value = "  -0.45%  "

$ws.Range("D40").Value = "'112.19"
$ws.Range("E40").Value = "  +0.68%  "

$ws.Range("D41").Value = "'2.026"
$ws.Range("E41").Value = "  -2.00%  "

$ws.Range("D42").Value = "'0.4443"
$ws.Range("E42").Value = "  +1.62%  "

$ws.Range("D43").Value = "'0.8651"
$ws.Range("E43").Value = "  -1.79%  "

$ws.Range("D44").Value = "'71.04"
$ws.Range("E44").Value = "  +5.45%  "

$ws.Range("D45").Value = "'5.844"
$ws.Range("E45").Value = "  -1.94%  "

$ws.Range("D46").Value = "'0.9997"
$ws.Range("E46").Value = "  -0.06%  "

$ws.Range("D47").Value = "'7.295"
$ws.Range("E47").Value = "  +0.21%  "

$ws.Range("D48").Value = "'48.71"
$ws.Range("E48").Value = "  +1.13%  "

$ws.Range("B49").Value = "Algorand"
$ws.Range("C49").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D49").Value = "'0.1236"
$ws.Range("E49").Value = "  -0.08%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'9.120"
$ws.Range("E50").Value = "  -2.16%  "

$ws.Range("D51").Value = "'0.2539"
$ws.Range("E51").Value = "  +3.49%  "
